$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in column F, row 1, matching the style
# of the existing header row (column E).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Timestamps for each data row (2-51), taken from the source commit.
$timestamps = @(
    "2021-10-05 10:51:11.900510",
    "2021-10-05 10:51:11.900523",
    "2021-10-05 10:51:11.900528",
    "2021-10-05 10:51:11.900532",
    "2021-10-05 10:51:11.900535",
    "2021-10-05 10:51:11.900539",
    "2021-10-05 10:51:11.900542",
    "2021-10-05 10:51:11.900545",
    "2021-10-05 10:51:11.900548",
    "2021-10-05 10:51:11.900552",
    "2021-10-05 10:51:11.900555",
    "2021-10-05 10:51:11.900558",
    "2021-10-05 10:51:11.900561",
    "2021-10-05 10:51:11.900564",
    "2021-10-05 10:51:11.900567",
    "2021-10-05 10:51:11.900570",
    "2021-10-05 10:51:11.900574",
    "2021-10-05 10:51:11.900577",
    "2021-10-05 10:51:11.900580",
    "2021-10-05 10:51:11.900583",
    "2021-10-05 10:51:11.900586",
    "2021-10-05 10:51:11.900589",
    "2021-10-05 10:51:11.900593",
    "2021-10-05 10:51:11.900596",
    "2021-10-05 10:51:11.900599",
    "2021-10-05 10:51:11.900602",
    "2021-10-05 10:51:11.900606",
    "2021-10-05 10:51:11.900609",
    "2021-10-05 10:51:11.900612",
    "2021-10-05 10:51:11.900615",
    "2021-10-05 10:51:11.900618",
    "2021-10-05 10:51:11.900621",
    "2021-10-05 10:51:11.900625",
    "2021-10-05 10:51:11.900628",
    "2021-10-05 10:51:11.900631",
    "2021-10-05 10:51:11.900634",
    "2021-10-05 10:51:11.900637",
    "2021-10-05 10:51:11.900641",
    "2021-10-05 10:51:11.900644",
    "2021-10-05 10:51:11.900647",
    "2021-10-05 10:51:11.900650",
    "2021-10-05 10:51:11.900654",
    "2021-10-05 10:51:11.900657",
    "2021-10-05 10:51:11.900660",
    "2021-10-05 10:51:11.900663",
    "2021-10-05 10:51:11.900666",
    "2021-10-05 10:51:11.900669",
    "2021-10-05 10:51:11.900672",
    "2021-10-05 10:51:11.900675",
    "2021-10-05 10:51:11.900678"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

Write-Host "Applied time_taken column to $($timestamps.Length) rows"
